$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

Set-TextValue 'D2' '27.243.08'
Set-TextValue 'E2' '  +0.17%  '
Set-TextValue 'D3' '1.906.21'
Set-TextValue 'E3' '  +0.70%  '
Set-TextValue 'E4' '  -0.13%  '
Set-TextValue 'D5' '306.35'
Set-TextValue 'E5' '  -0.17%  '
Set-TextValue 'E6' '  -0.05%  '
Set-TextValue 'E7' '  +2.95%  '
Set-TextValue 'D8' '0.3817'
Set-TextValue 'E8' '  +1.76%  '
Set-TextValue 'D9' '0.07293'
Set-TextValue 'E9' '  +0.49%  '
Set-TextValue 'E10' '  +5.11%  '
Set-TextValue 'D11' '0.9055'
Set-TextValue 'E11' '  +0.90%  '
Set-TextValue 'D12' '0.08195'
Set-TextValue 'E12' '  +0.42%  '
Set-TextValue 'D13' '95.80'
Set-TextValue 'E13' '  -0.93%  '
Set-TextValue 'D14' '5.356'
Set-TextValue 'E14' '  +1.71%  '
Set-TextValue 'E15' '  -0.09%  '
Set-TextValue 'D16' '14.87'
Set-TextValue 'E16' '  +2.39%  '
Set-TextValue 'D17' '0.000008676'
Set-TextValue 'E17' '  +1.13%  '
Set-TextValue 'E18' '  -0.08%  '
Set-TextValue 'D19' '27.259.11'
Set-TextValue 'D20' '5.050'
Set-TextValue 'E20' '  -0.54%  '
Set-TextValue 'D21' '1.085.88'
Set-TextValue 'E21' '  -42.63%  '
Set-TextValue 'D22' '10.79'
Set-TextValue 'E22' '  +1.15%  '
Set-TextValue 'D23' '6.525'
Set-TextValue 'E23' '  +2.02%  '
Set-TextValue 'D24' '149.02'
Set-TextValue 'E24' '  +1.08%  '
Set-TextValue 'D25' '2.305'
Set-TextValue 'E25' '  +0.89%  '
Set-TextValue 'E26' '  +1.29%  '
Set-TextValue 'D27' '1.748'
Set-TextValue 'E27' '  +0.23%  '
Set-TextValue 'D28' '116.77'
Set-TextValue 'E28' '  +1.63%  '
Set-TextValue 'D29' '4.834'
Set-TextValue 'E29' '  +0.99%  '
Set-TextValue 'D30' '4.721'
Set-TextValue 'E30' '  -3.84%  '
Set-TextValue 'D31' '0.09219'
Set-TextValue 'E31' '  -0.01%  '
Set-TextValue 'D32' '0.8295'
Set-TextValue 'E32' '  +5.02%  '
Set-TextValue 'D33' '0.05087'
Set-TextValue 'E33' '  +1.09%  '
Set-TextValue 'D34' '1.218'
Set-TextValue 'E34' '  +0.22%  '
Set-TextValue 'D35' '2.994'
Set-TextValue 'E35' '  +1.51%  '
Set-TextValue 'E36' '  -3.44%  '
Set-TextValue 'D37' '2.671'
Set-TextValue 'E37' '  +4.29%  '
Set-TextValue 'D38' '0.5942'
Set-TextValue 'E38' '  +5.39%  '
Set-TextValue 'E39' '  +1.11%  '
Set-TextValue 'D40' '1.079'
Set-TextValue 'E40' '  +0.67%  '
Set-TextValue 'D41' '9.371'
Set-TextValue 'E41' '  +5.06%  '
Set-TextValue 'D42' '6.678'
Set-TextValue 'E42' '  +2.35%  '
Set-TextValue 'D43' '116.88'
Set-TextValue 'E43' '  +1.52%  '
Set-TextValue 'D44' '0.5144'
Set-TextValue 'E44' '  +5.92%  '
Set-TextValue 'E45' '  +1.19%  '
Set-TextValue 'D46' '10.23'
Set-TextValue 'E46' '  +1.93%  '
Set-TextValue 'E47' '  -0.03%  '
Set-TextValue 'E48' '  +1.90%  '
Set-TextValue 'D49' '38.46'
Set-TextValue 'E49' '  +1.03%  '
Set-TextValue 'D50' '0.06148'
Set-TextValue 'E50' '  +3.48%  '
Set-TextValue 'D51' '63.51'
Set-TextValue 'E51' '  +0.41%  '
